$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "username" key row
$ws.Range("A7").Value = "username"

# Fix the description text for "selectedCategory" (B4): "beinhaltet" -> "enthält"
$ws.Range("B4").Value = "enthält die durch den Nutzer (nicht Gegner) ausgewählte Kategorie"

# Add "password" key row
$ws.Range("A8").Value = "password"

# Fill in the value descriptions for the new rows
$ws.Range("B7").Value = "enthält den nutzernamen des angemeldeten benutzers"
$ws.Range("B8").Value = "enthält das passwort des angemeldeten benutzers"

# Resize the table to include the two new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A3:B8"))

# Update selection to match the post-edit state (active cell moved to B9)
$ws.Range("B9").Select()

$wb.Save()
